# Appendix Table S2 update:
#  - update table numbers to reflect the manuscript renumbering
#  - change the Greek alpha glyph in "Edible Chlorophyll α" to an
#    italicized Latin "a" ("Edible Chlorophyll a")
#  - remove the leftover picture that used to sit next to the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the embedded picture/drawing object from the worksheet.
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Replace "Edible Chlorophyll α" with an italicized "Edible Chlorophyll a"
# in every cell that references it.
$cells = @("C4", "B6", "B7")
foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $cell.Value = "Edible Chlorophyll a"
    # Italicize just the trailing "a" (position 20, length 1) to match
    # the "Edible Chlorophyll " + italic "a" run split.
    $cell.Characters(20, 1).Font.Italic = $true
}
